$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "42.391.48"
$ws.Range("E2").Value = "  +0.54%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.272.92"
$ws.Range("E3").Value = "  -0.35%  "

$ws.Range("E4").Value = "  +0.08%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "306.70"
$ws.Range("E5").Value = "  +0.45%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "98.29"
$ws.Range("E6").Value = "  +3.09%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.528"
$ws.Range("E7").Value = "  -0.70%  "

$ws.Range("E8").Value = "  +0.07%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.495"
$ws.Range("E9").Value = "  +0.23%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "35.52"
$ws.Range("E10").Value = "  -0.43%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0791"
$ws.Range("E11").Value = "  -1.57%  "

$ws.Range("E12").Value = "  -0.11%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.94"
$ws.Range("E13").Value = "  +3.49%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.625.79"
$ws.Range("E14").Value = "  -0.29%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "14.75"
$ws.Range("E15").Value = "  +1.72%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.290.85"
$ws.Range("E16").Value = "  +0.85%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.795"
$ws.Range("E17").Value = "  -0.56%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "42.249.67"
$ws.Range("E18").Value = "  +0.42%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.51"
$ws.Range("E19").Value = "  -1.98%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0907"
$ws.Range("E20").Value = "  -1.13%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.04"
$ws.Range("E21").Value = "  +0.20%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "68.17"
$ws.Range("E22").Value = "  +0.01%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "238.61"
$ws.Range("E23").Value = "  -2.09%  "

$ws.Range("E24").Value = "  +2.10%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.58"
$ws.Range("E25").Value = "  -0.83%  "

$ws.Range("E26").Value = "  -0.03%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "23.65"
$ws.Range("E27").Value = "  -1.79%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "37.94"
$ws.Range("E28").Value = "  +4.89%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.54"
$ws.Range("E29").Value = "  -1.83%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.12"
$ws.Range("E30").Value = "  +0.92%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "161.75"
$ws.Range("E31").Value = "  -0.03%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.26"
$ws.Range("E32").Value = "  -2.05%  "

$ws.Range("E33").Value = "  +0.13%  "

$ws.Range("E34").Value = "  +2.62%  "

$ws.Range("E35").Value = "  -1.31%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "17.71"
$ws.Range("E36").Value = "  +2.96%  "

$ws.Range("E37").Value = "  -0.63%  "

$ws.Range("E38").Value = "  -3.82%  "

$ws.Range("E39").Value = "  +0.35%  "

$ws.Range("E40").Value = "  -1.41%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.09"
$ws.Range("E41").Value = "  -2.36%  "

$ws.Range("E42").Value = "  +2.89%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.947.41"
$ws.Range("E43").Value = "  -3.22%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0283"
$ws.Range("E44").Value = "  -0.58%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "18.85"
$ws.Range("E45").Value = "  -3.59%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "9.94"
$ws.Range("E46").Value = "  -2.88%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.91"
$ws.Range("E47").Value = "  -3.90%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "53.83"
$ws.Range("E48").Value = "  +0.36%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "92.43"
$ws.Range("E49").Value = "  -0.33%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "71.95"
$ws.Range("E50").Value = "  -0.62%  "

$ws.Range("E51").Value = "  -1.91%  "
